$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the three new protocol rows (34-36) describing the "launch spell"
# messages and the PA-sync message. Source/Destination (A/B) reuse the
# existing Client/Serveur shared strings, so they can be written any time;
# the Paramètres/Description pairs go in first, then row 36 end-to-end, and
# the "CL" message code (shared by rows 34 & 35) is typed in last.
$ws.Range("A34").Value = "Client"
$ws.Range("B34").Value = "Serveur"
$ws.Range("D34").Value = "idSpell;cellX;cellY"
$ws.Range("E34").Value = "Le client indique qu'il veut lancer un sort."

$ws.Range("A35").Value = "Serveur"
$ws.Range("B35").Value = "Client"
$ws.Range("D35").Value = "idPerso;idSpell;cellX;cellY"
$ws.Range("E35").Value = "Le serveur indique qu'un personnage lance un sort."

$ws.Range("A36").Value = "Serveur"
$ws.Range("B36").Value = "Client"
$ws.Range("C36").Value = "Ca"
$ws.Range("D36").Value = "idPerso;nbPA"
$ws.Range("E36").Value = "Le serveur envoi l'information du nombre de PA restant du personnage actif (synchronisation quand il y a connexion en cours de combat)"

$ws.Range("C34").Value = "CL"
$ws.Range("C35").Value = "CL"

# Rows 35/36 now wrap across more lines than before, so their row heights
# grow to fit the new text (rows 2-33 follow the same ht="30"/"45" pattern).
$ws.Rows.Item(35).RowHeight = 30
$ws.Rows.Item(36).RowHeight = 45

# Scroll the frozen pane down one row and move the active selection to
# reflect where the author ended up working.
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("C40").Select() | Out-Null
